# Applies the cell-value updates described by the commit diff.
# Sheet order: 1=exhibitions(展览) 2=shows(演出) 3=local-life(本地生活) 4=all-types(全部类型, union of the above)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# --- Sheet 1: 展览 ---
$ws1.Range("F2").Value = 316
$ws1.Range("F3").Value = 518
$ws1.Range("F4").Value = 509
$ws1.Range("F5").Value = 9006
$ws1.Range("F6").Value = 9006
$ws1.Range("F8").Value = 11609
$ws1.Range("F14").Value = 130
$ws1.Range("F16").Value = 306
$ws1.Range("F17").Value = 33
$ws1.Range("F18").Value = 102
$ws1.Range("F20").Value = 433
$ws1.Range("F21").Value = 1978
$ws1.Range("F22").Value = 764
$ws1.Range("F23").Value = 708
$ws1.Range("F24").Value = 376
$ws1.Range("F25").Value = 10
$ws1.Range("F26").Value = 303
$ws1.Range("F29").Value = 1426
$ws1.Range("F30").Value = 34
$ws1.Range("F33").Value = 50
$ws1.Range("F34").Value = 2
$ws1.Range("F35").Value = 1431
$ws1.Range("F37").Value = 485
$ws1.Range("F38").Value = 328
$ws1.Range("F39").Value = 177
$ws1.Range("F40").Value = 368
$ws1.Range("F41").Value = 52
$ws1.Range("F42").Value = 152
$ws1.Range("F43").Value = 557
$ws1.Range("F44").Value = 420
$ws1.Range("F45").Value = 139
$ws1.Range("F46").Value = 844
$ws1.Range("F49").Value = 216
$ws1.Range("F50").Value = 194

# --- Sheet 2: 演出 ---
$ws2.Range("F8").Value = 67
$ws2.Range("F17").Value = 69
$ws2.Range("F24").Value = 56
$ws2.Range("F25").Value = 97
$ws2.Range("F26").Value = 17
$ws2.Range("F27").Value = 412

# --- Sheet 3: 本地生活 ---
$ws3.Range("F3").Value = 2885
$ws3.Range("F4").Value = 359
$ws3.Range("F5").Value = 228
$ws3.Range("F6").Value = 94
$ws3.Range("G3").Value = "已售罄"

# --- Sheet 4: 全部类型 ---
$ws4.Range("F2").Value = 518
$ws4.Range("F4").Value = 228
$ws4.Range("F6").Value = 509
$ws4.Range("F7").Value = 9006
$ws4.Range("F8").Value = 9006
$ws4.Range("F10").Value = 11609
$ws4.Range("F15").Value = 130
$ws4.Range("F16").Value = 306
$ws4.Range("F17").Value = 102
$ws4.Range("F19").Value = 1978
$ws4.Range("F20").Value = 764
$ws4.Range("F22").Value = 376
$ws4.Range("F23").Value = 303
$ws4.Range("F26").Value = 67
$ws4.Range("F28").Value = 1426
$ws4.Range("F29").Value = 34
$ws4.Range("F33").Value = 50
$ws4.Range("F35").Value = 1431
$ws4.Range("F36").Value = 328
$ws4.Range("B37").Value = "'2024-10-30"
$ws4.Range("C37").Value = "北京·majiko中国巡演-2024"
$ws4.Range("D37").Value = "奥园西路1号院5号楼1层2-104 福浪Live House"
$ws4.Range("E37").Value = "2024.10.30 20:00-10.30 21:40"
$ws4.Range("F37").Value = 79
$ws4.Range("G37").Value = 480
$ws4.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=92300"
$ws4.Range("I37").Value = "//i0.hdslb.com/bfs/openplatform/202409/sxKGqhKo1726038853275.jpeg"
$ws4.Range("B38").Value = "'2024-11-02"
$ws4.Range("C38").Value = "北京·明日方舟only同人 2.0 京台夕照"
$ws4.Range("D38").Value = "永外高庄138号 北京大红门国际会展中心"
$ws4.Range("E38").Value = "2024.11.02 09:00-11.02 17:00"
$ws4.Range("F38").Value = 368
$ws4.Range("G38").Value = 89
$ws4.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=90479"
$ws4.Range("I38").Value = "//i2.hdslb.com/bfs/openplatform/202408/TMycI1on1723101208256.jpeg"
$ws4.Range("B39").Value = "'2024-11-03"
$ws4.Range("C39").Value = "北京·从“梁祝”到“吉普赛之歌”——小提琴王子刘霄经典名曲音乐会"
$ws4.Range("D39").Value = "中关村南大街33号国家图书馆内 国图艺术中心"
$ws4.Range("E39").Value = "2024.11.03 19:30-11.03 21:00"
$ws4.Range("F39").Value = 1
$ws4.Range("G39").Value = 126
$ws4.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=91379"
$ws4.Range("I39").Value = "//i0.hdslb.com/bfs/openplatform/202408/Lro1XCdG1724743023171.jpeg"
$ws4.Range("B40").Value = "'2024-11-09"
$ws4.Range("C40").Value = "北京·漫视界IIS动漫游戏盛典02"
$ws4.Range("D40").Value = "黑庄户路8号 北京音乐产业园"
$ws4.Range("E40").Value = "2024.11.09 09:30-11.10 17:00"
$ws4.Range("F40").Value = 152
$ws4.Range("G40").Value = 76
$ws4.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=90029"
$ws4.Range("I40").Value = "//i0.hdslb.com/bfs/openplatform/202407/QfjWfDU71721980061932.jpeg"
$ws4.Range("F42").Value = 420
$ws4.Range("F43").Value = 139
$ws4.Range("F45").Value = 56
$ws4.Range("F46").Value = 412
$ws4.Range("F49").Value = 216
$ws4.Range("F50").Value = 194

